# Auto commit at 2025-12-02 8:00:30.79
# Update Metrics figures, refresh dependent "today" formulas, and move the
# active-sheet/selection state from "ndzsrqs" to "today".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Metrics sheet: update the raw figures in column B (rows 2-13).
# ---------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 12173.83
$wsMetrics.Range("B3").Value  = 10419.36
$wsMetrics.Range("B4").Value  = 3656.01
$wsMetrics.Range("B5").Value  = 494
$wsMetrics.Range("B6").Value  = 5214880.94
$wsMetrics.Range("B7").Value  = 4410772.32
$wsMetrics.Range("B8").Value  = 1535612.89
$wsMetrics.Range("B9").Value  = 203201
$wsMetrics.Range("B10").Value = 33680261.93
$wsMetrics.Range("B11").Value = 31686047.48
$wsMetrics.Range("B12").Value = 11817334.93
$wsMetrics.Range("B13").Value = 1300831

# Recalculate so dependent formulas (on the "today" sheet) pick up the
# new values right away.
$excel.Calculate()

# Move the Metrics selection (no tab activation here; "today" ends up
# being the active sheet below).
[void]$wsMetrics.Range("D21").Select()

# ---------------------------------------------------------------------
# 2) "today" sheet becomes the active / selected sheet, with a new
#    selected cell.
# ---------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Activate()
[void]$wsToday.Range("G20").Select()

# ---------------------------------------------------------------------
# 3) "ndzsrqs" sheet loses the tab-selected flag, but keeps its own
#    last selection (L20) untouched.
# ---------------------------------------------------------------------
$wsNdzsrqs = $wb.Worksheets.Item("ndzsrqs")
[void]$wsNdzsrqs.Range("L20").Select()

# Re-activate "today" last, so it is the sheet that ends up marked as
# active/tabSelected in the saved workbook.
$wsToday.Activate()

$excel.Calculate()
